# Insert a new price-record row at row 159 (pushes existing rows 159..274
# down to 160..275) and populate it with the new observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(159).Insert()

$ws.Cells.Item(159, 1).Value  = 3
$ws.Cells.Item(159, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(159, 3).Value  = "Coquimbo"
$ws.Cells.Item(159, 4).Value  = 44574
$ws.Cells.Item(159, 5).Value  = 5
$ws.Cells.Item(159, 6).Value  = 100114013
$ws.Cells.Item(159, 7).Value  = "Zanahoria"
$ws.Cells.Item(159, 8).Value  = "Sin especificar"
$ws.Cells.Item(159, 9).Value  = "Primera"
$ws.Cells.Item(159, 10).Value = 225
$ws.Cells.Item(159, 11).Value = 7000
$ws.Cells.Item(159, 12).Value = 7500
$ws.Cells.Item(159, 13).Value = 7278
$ws.Cells.Item(159, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(159, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(159, 16).Value = 364
$ws.Cells.Item(159, 17).Value = 20
$ws.Cells.Item(159, 18).Value = "Hortaliza"
